$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 (FEED / External Devices): the adjoining "**Cons:**" note that used
# to follow this text in the table was removed, so the now-dangling
# "--- --- --- ---" separator at the end of this cell is trimmed off too.
$ws.Range("D6").Value = "**Pros:**  Can connect to external sensors; uses local area network "

# Row 7 (Big Fin / Data Interface): drop the stray "**Cons:**  " placeholder
# cell entirely.
[void]$ws.Range("C7").ClearContents()

# Update the saved view state: scroll up a couple of rows and leave the
# active selection on G4.
$win = $excel.ActiveWindow
[void]$ws.Range("G4").Select()
$win.ScrollRow = 3
$win.ScrollColumn = 1
